$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Enterprises density (per 1000 people) -> row 13
$ws.Range("B13").Value = "'37.86"
$ws.Range("C13").Value = "'4.13"
$ws.Range("D13").Value = "'41.99"

# Employment (% of total) -> row 14
$ws.Range("B14").Value = "'24.86"
$ws.Range("C14").Value = "'51.54"
$ws.Range("D14").Value = "'76.39"

# Enterprises (% of total) -> row 16
$ws.Range("B16").Value = "'89.97"
$ws.Range("C16").Value = "'9.81"
$ws.Range("D16").Value = "'99.77"

# Value added to the economy (% of total) -> row 20
$ws.Range("B20").Value = "'13.48"
$ws.Range("C20").Value = "'55.49"
$ws.Range("D20").Value = "'68.98"
